# Update the test credentials used on the LoginCredentials sheet:
# the vsuser/Govirtual~1 pair is replaced with shopfloor1/shopfloor1*1
# for both the QA and RC environment rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginCredentials")

$ws.Range("B2").Value = "shopfloor1"
$ws.Range("C2").Value = "shopfloor1*1"
$ws.Range("B3").Value = "shopfloor1"
$ws.Range("C3").Value = "shopfloor1*1"

# Column C needs to widen slightly to fit the new, one-character-longer
# password value ("shopfloor1*1" vs "Govirtual~1").
$ws.Range("C1:C3").ColumnWidth = 11.666666666666668

# Leave the same cell selection state the sheet was saved with.
$ws.Range("B3:C3").Select()
